# Update cryptos list data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.646.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.193.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.25%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.82%  '

$ws.Range("E7").Value = '  +1.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.692'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.46%  '

$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.188.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.26%  '

$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.177'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.73%  '

$ws.Range("E13").Value = '  -4.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.480.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.773.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.22%  '

$ws.Range("E16").Value = '  -4.34%  '

$ws.Range("E17").Value = '  -3.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.187.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.59%  '

$ws.Range("E19").Value = '  +5.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '444.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000187'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +35.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.77%  '

$ws.Range("E24").Value = '  -4.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.25%  '

$ws.Range("E26").Value = '  -3.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '75.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("E30").Value = '  -7.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +33.09%  '

$ws.Range("E33").Value = '  -4.92%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '535.11'
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.26%  '

$ws.Range("E36").Value = '  -5.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.29%  '

$ws.Range("E40").Value = '  -9.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("E43").Value = '  -6.70%  '

$ws.Range("E44").Value = '  -6.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '44.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '172.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.19%  '

$ws.Range("E48").Value = '  -6.04%  '

$ws.Range("E49").Value = '  -5.47%  '

$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.83%  '

$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.614'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.65%  '
